$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 291
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 22)  # column V = 22
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val * 100
    }
}
